# [Fonds de solidarite] Add 2020-08-25 data
# Updates "nombre_aides" (C) and "montant_total" (D) for the rows whose
# regional/classe_effectif bucket received revised counts in the
# 2020-08-25 refresh of the Fonds de solidarite volet 2 dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; C = "190"; D = "448016.00" },
    @{ Row = 3; C = "1010"; D = "3223895.33" },
    @{ Row = 4; C = "420"; D = "1735698.25" },
    @{ Row = 8; C = "42"; D = "84000.00" },
    @{ Row = 9; C = "58"; D = "149597.64" },
    @{ Row = 10; C = "358"; D = "1274188.71" },
    @{ Row = 33; C = "109"; D = "317173.00" },
    @{ Row = 34; C = "570"; D = "1877479.47" },
    @{ Row = 35; C = "227"; D = "1144788.11" },
    @{ Row = 38; C = "22"; D = "48200.00" },
    @{ Row = 46; C = "85"; D = "374474.61" },
    @{ Row = 50; C = "14"; D = "31850.00" },
    @{ Row = 51; C = "103"; D = "295768.17" },
    @{ Row = 52; C = "595"; D = "2098936.52" },
    @{ Row = 53; C = "263"; D = "1152878.76" },
    @{ Row = 57; C = "712"; D = "1821318.62" },
    @{ Row = 58; C = "3525"; D = "11041660.38" },
    @{ Row = 59; C = "1824"; D = "7601634.94" },
    @{ Row = 60; C = "634"; D = "3120291.45" },
    @{ Row = 61; C = "128"; D = "885123.00" },
    @{ Row = 63; C = "285"; D = "676120.58" },
    @{ Row = 82; C = "229"; D = "583326.09" },
    @{ Row = 83; C = "889"; D = "2850012.26" }
)

foreach ($u in $updates) {
    $cCell = $ws.Cells.Item($u.Row, 3)   # column C: nombre_aides
    $cCell.NumberFormat = "@"
    $cCell.Value = $u.C

    $dCell = $ws.Cells.Item($u.Row, 4)   # column D: montant_total
    $dCell.NumberFormat = "@"
    $dCell.Value = $u.D
}
